$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-30 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-31 Friday", 2) | Out-Null
$d.Content.Find.Execute("505÷8=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "235÷8=29, 3", 2) | Out-Null
$d.Content.Find.Execute("461÷3=153, 2", $true, $false, $false, $false, $false, $true, 1, $false, "910÷9=101, 1", 2) | Out-Null
$d.Content.Find.Execute("124÷4=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "822÷4=205, 2", 2) | Out-Null
$d.Content.Find.Execute("274÷4=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "432÷2=216, 0", 2) | Out-Null
$d.Content.Find.Execute("512÷8=64, 0", $true, $false, $false, $false, $false, $true, 1, $false, "304÷3=101, 1", 2) | Out-Null
$d.Content.Find.Execute("377÷6=62, 5", $true, $false, $false, $false, $false, $true, 1, $false, "122÷3=40, 2", 2) | Out-Null
$d.Content.Find.Execute("565÷5=113, 0", $true, $false, $false, $false, $false, $true, 1, $false, "940÷8=117, 4", 2) | Out-Null
$d.Content.Find.Execute("480÷8=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "206÷4=51, 2", 2) | Out-Null
$d.Content.Find.Execute("933÷4=233, 1", $true, $false, $false, $false, $false, $true, 1, $false, "203÷8=25, 3", 2) | Out-Null
$d.Content.Find.Execute("730÷8=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "158÷2=79, 0", 2) | Out-Null
$d.Content.Find.Execute("737÷5=147, 2", $true, $false, $false, $false, $false, $true, 1, $false, "773÷8=96, 5", 2) | Out-Null
$d.Content.Find.Execute("279÷3=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "854÷2=427, 0", 2) | Out-Null
$d.Content.Find.Execute("302÷4=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "724÷3=241, 1", 2) | Out-Null
$d.Content.Find.Execute("386÷5=77, 1", $true, $false, $false, $false, $false, $true, 1, $false, "347÷6=57, 5", 2) | Out-Null
$d.Content.Find.Execute("435÷2=217, 1", $true, $false, $false, $false, $false, $true, 1, $false, "990÷9=110, 0", 2) | Out-Null
$d.Content.Find.Execute("670÷5=134, 0", $true, $false, $false, $false, $false, $true, 1, $false, "209÷2=104, 1", 2) | Out-Null
$d.Content.Find.Execute("812÷5=162, 2", $true, $false, $false, $false, $false, $true, 1, $false, "954÷3=318, 0", 2) | Out-Null
$d.Content.Find.Execute("156÷5=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "496÷7=70, 6", 2) | Out-Null
$d.Content.Find.Execute("392÷8=49, 0", $true, $false, $false, $false, $false, $true, 1, $false, "204÷9=22, 6", 2) | Out-Null
$d.Content.Find.Execute("389÷3=129, 2", $true, $false, $false, $false, $false, $true, 1, $false, "781÷2=390, 1", 2) | Out-Null
$d.Content.Find.Execute("571÷8=71, 3", $true, $false, $false, $false, $false, $true, 1, $false, "693÷5=138, 3", 2) | Out-Null
$d.Content.Find.Execute("832÷9=92, 4", $true, $false, $false, $false, $false, $true, 1, $false, "445÷6=74, 1", 2) | Out-Null
$d.Content.Find.Execute("651÷9=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "992÷3=330, 2", 2) | Out-Null
$d.Content.Find.Execute("234÷4=58, 2", $true, $false, $false, $false, $false, $true, 1, $false, "398÷8=49, 6", 2) | Out-Null
$d.Content.Find.Execute("543÷4=135, 3", $true, $false, $false, $false, $false, $true, 1, $false, "295÷3=98, 1", 2) | Out-Null
